$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new row 23 with October 2025 stats, copying the date style from A22
$ws.Range("A23").Value = 45931
$ws.Range("A23").NumberFormat = $ws.Range("A22").NumberFormat

$ws.Range("B23").Value = 6303
$ws.Range("C23").Value = 996
$ws.Range("D23").Value = 5866825
$ws.Range("E23").Value = 930.798825955894
$ws.Range("F23").Value = 8.150308853809207
$ws.Range("G23").Value = 3.642039542143594
$ws.Range("H23").Value = 25.71104282423948
